$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5539.9085,
    5343.505999999999,
    5250.640500000001,
    5194.129499999999,
    5143.397,
    5188.8305,
    5628.136500000001,
    6549.074000000001,
    7809.641,
    12977.8705,
    14879.494,
    15334.284,
    15538.4565,
    15473.9585,
    15459.3285,
    15733.5675,
    15999.228,
    15687.07,
    15318.933,
    13434.6415,
    11609.815,
    9103.451000000001,
    6290.4415,
    5217.8665
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
